$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D6","D7","D9","D10","D11","D12","D14","D15","D16","D17","D20","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.690.35"
$ws.Range("E2").Value = "  +5.08%  "
$ws.Range("D3").Value = "2.230.87"
$ws.Range("E3").Value = "  +3.39%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "227.65"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("D7").Value = "60.57"
$ws.Range("E7").Value = "  -4.60%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "0.401"
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").Value = "58.11"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").Value = "0.0871"
$ws.Range("E11").Value = "  +3.04%  "
$ws.Range("D12").Value = "0.103"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "2.571.97"
$ws.Range("E13").Value = "  +3.58%  "
$ws.Range("D14").Value = "15.63"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").Value = "21.35"
$ws.Range("E15").Value = "  -2.48%  "
$ws.Range("D16").Value = "0.794"
$ws.Range("E16").Value = "  -1.51%  "
$ws.Range("D17").Value = "5.53"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "2.250.03"
$ws.Range("E18").Value = "  +3.90%  "
$ws.Range("D19").Value = "41.629.51"
$ws.Range("E19").Value = "  +5.12%  "
$ws.Range("D20").Value = "72.51"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").Value = "0.0₃0883"
$ws.Range("E21").Value = "  +4.71%  "
$ws.Range("D22").Value = "6.02"
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("D23").Value = "246.32"
$ws.Range("E23").Value = "  +6.62%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("D26").Value = "2.32"
$ws.Range("E26").Value = "  -1.95%  "
$ws.Range("D27").Value = "9.49"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").Value = "167.29"
$ws.Range("E28").Value = "  -3.19%  "
$ws.Range("D29").Value = "0.140"
$ws.Range("E29").Value = "  +1.57%  "
$ws.Range("D30").Value = "19.87"
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("D31").Value = "1.41"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").Value = "2.65"
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("D33").Value = "0.122"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D34").Value = "5.04"
$ws.Range("E34").Value = "  +8.29%  "
$ws.Range("D35").Value = "4.64"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").Value = "0.0619"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").Value = "6.59"
$ws.Range("E37").Value = "  -5.05%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "3.68"
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("D39").Value = "2.36"
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("B40").Value = "BinanceUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D40").Value = "1.01"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("B41").Value = "TerraClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D41").Value = "0.000237"
$ws.Range("E41").Value = "  +28.94%  "
$ws.Range("D42").Value = "4.88"
$ws.Range("E42").Value = "  -4.95%  "
$ws.Range("D43").Value = "0.0234"
$ws.Range("E43").Value = "  +3.45%  "
$ws.Range("D44").Value = "8.65"
$ws.Range("E44").Value = "  +11.66%  "
$ws.Range("D45").Value = "0.0969"
$ws.Range("E45").Value = "  +5.46%  "
$ws.Range("D46").Value = "98.90"
$ws.Range("E46").Value = "  -3.20%  "
$ws.Range("D47").Value = "1.465.78"
$ws.Range("E47").Value = "  -3.01%  "
$ws.Range("D48").Value = "1.18"
$ws.Range("E48").Value = "  -3.32%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "2.78"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "16.18"
$ws.Range("E50").Value = "  -7.00%  "
$ws.Range("D51").Value = "1.08"
$ws.Range("E51").Value = "  -1.05%  "
